$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.352.56"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "1.840.99"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.99"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6262"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -0.71%  "
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.76"
$ws.Range("E10").Value = "  +1.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07711"
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("D12").Value = "1.842.13"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.964"
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6742"
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001025"
$ws.Range("E15").Value = "  -1.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.76"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.238"
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("D18").Value = "29.327.74"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "234.32"
$ws.Range("E19").Value = "  +2.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.32"
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.308"
$ws.Range("E22").Value = "  -2.71%  "
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "157.87"
$ws.Range("E24").Value = "  -0.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.481"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1345"
$ws.Range("E26").Value = "  -1.53%  "
$ws.Range("E27").Value = "  -1.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.07232"
$ws.Range("E28").Value = "  +10.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.472"
$ws.Range("E29").Value = "  +3.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.476"
$ws.Range("E30").Value = "  -0.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.063"
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.025"
$ws.Range("E32").Value = "  -1.50%  "
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.143"
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6943"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01834"
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.912"
$ws.Range("E38").Value = "  +2.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.810"
$ws.Range("E39").Value = "  -0.88%  "
$ws.Range("D40").Value = "1.232.02"
$ws.Range("E40").Value = "  -2.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9476"
$ws.Range("E41").Value = "  +3.35%  "
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("D43").Value = "1.995.31"
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.89"
$ws.Range("E44").Value = "  -0.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.17"
$ws.Range("E45").Value = "  -1.15%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.705"
$ws.Range("E46").Value = "  -1.63%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.934"
$ws.Range("E47").Value = "  -1.80%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000115"
$ws.Range("E48").Value = "  -2.55%  "
$ws.Range("E49").Value = "  -1.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.812"
$ws.Range("E50").Value = "  -1.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1126"
$ws.Range("E51").Value = "  -2.51%  "
